$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text change
$ws.Range("B1").Value = "Load (KN)"

# New data values (rows 2..19); rows 20..22 are cleared (table shrinks from 22 to 19 rows)
$data = @(
    @(10, 2333.65),
    @(74, 3267.765849889729),
    @(66, 3196.18),
    @(51, 1263),
    @(7, 1640),
    @(78, 2334.24),
    @(92, 2013),
    @(90, 3968),
    @(107, 2173),
    @(62, 1279),
    @(30, 1193.56),
    @(2, 1639.98),
    @(94, 3359.28),
    @(22, 1816.34),
    @(43, 1340),
    @(84, 3891.95),
    @(24, 1705.51),
    @(73, 2017.001651742916)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the now-unused trailing rows (20, 21, 22) so the sheet dimension shrinks to A1:B19
$ws.Range("A20:B22").Clear()
